$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Oil (O), Protein (P), and Oil+Protein (Q) lab-result values for
# reps 1 and 2 of the first 20 genotypes (test weight / lsmeans cleanup).
$oilProteinData = @{
    2 = @(21.94, 43.9, 65.84)
    3 = @(21.64, 44.82, 66.46)
    6 = @(22.98, 41.96, 64.94)
    7 = @(23.28, 40.8, 64.08)
    10 = @(22.03, 41.22, 63.25)
    11 = @(21.86, 43.04, 64.9)
    14 = @(20.48, 44.87, 65.35)
    15 = @(20.86, 44.9, 65.76)
    18 = @(22.6, 41.52, 64.12)
    19 = @(21.92, 43.39, 65.31)
    22 = @(21.74, 44.65, 66.39)
    23 = @(21.61, 45.58, 67.19)
    26 = @(22.71, 42.6, 65.31)
    27 = @(23.24, 42.79, 66.03)
    30 = @(22.48, 42.46, 64.94)
    31 = @(22.85, 43.19, 66.04)
    34 = @(21.76, 41.78, 63.54)
    35 = @(20.85, 43.2, 64.05)
    38 = @(21.19, 45.81, 67)
    39 = @(21.18, 46.94, 68.12)
    42 = @(22.02, 42.34, 64.36)
    43 = @(22.45, 40.95, 63.4)
    46 = @(20.89, 45.11, 66)
    47 = @(20.54, 46.5, 67.04)
    50 = @(21.76, 43.94, 65.7)
    51 = @(21.46, 43.91, 65.37)
    54 = @(21.6, 44.18, 65.78)
    55 = @(21.21, 44.35, 65.56)
    58 = @(21.69, 44.59, 66.28)
    59 = @(22.35, 44.54, 66.89)
    62 = @(23.07, 39, 62.07)
    63 = @(22.44, 39.17, 61.61)
    66 = @(22.91, 39.28, 62.19)
    67 = @(23.15, 38.52, 61.67)
    70 = @(23.73, 39.36, 63.09)
    71 = @(24.51, 38.37, 62.88)
    74 = @(19.31, 49.13, 68.44)
    75 = @(18.52, 51.32, 69.84)
    78 = @(20.85, 44.78, 65.63)
    79 = @(20.27, 44.67, 64.94)
}

foreach ($row in $oilProteinData.Keys) {
    $values = $oilProteinData[$row]
    $ws.Cells.Item($row, 15).Value = $values[0]  # Column O - Oil
    $ws.Cells.Item($row, 16).Value = $values[1]  # Column P - Protein
    $ws.Cells.Item($row, 17).Value = $values[2]  # Column Q - Oil + Protein
}
